$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Stage the literal text "true" in a scratch cell via a formula (so it is a
# real String value, not user-typed input) and copy/paste-values it into the
# target cells. This avoids Excel's automatic Boolean coercion that a direct
# Range.Value = "true" assignment would trigger, while keeping each cell on
# the sheet's existing body style (no quote-prefixed text, no new style).
$ws.Range("D1").Formula = "=""true"""
$ws.Range("D1").Copy()

# "Experimental" flips from false to true
$ws.Range("B7").PasteSpecial(-4163)

# "Case Sensitive" row gains a "true" value as well
$ws.Range("B14").PasteSpecial(-4163)

$ws.Range("D1").ClearContents()
$excel.CutCopyMode = $false

# Publication date is refreshed to the new build timestamp
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
